$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Delete column D -> shifts E:M left to D:L, and updates col widths/dimension automatically
$ws.Columns("D").Delete()

# Seed new column M's formatting (styles) by copying column L's formats (used range only)
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match column M width to the repeating-pattern width (29, same as col E/I)
$ws.Columns("M").ColumnWidth = $ws.Columns("E").ColumnWidth

# Row 8: financial period headers (shifted + newest quarter appended)
$ws.Range("D8").Value = "9 ماهه منتهی به 1399/09"
$ws.Range("E8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("F8").Value = "3 ماهه منتهی به 1400/03"
$ws.Range("G8").Value = "6 ماهه منتهی به 1400/06"
$ws.Range("H8").Value = "9 ماهه منتهی به 1400/09"
$ws.Range("I8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("J8").Value = "3 ماهه منتهی به 1401/03"
$ws.Range("K8").Value = "6 ماهه منتهی به 1401/06"
$ws.Range("L8").Value = "9 ماهه منتهی به 1401/09"
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# Row 9: publish dates (shifted + newest publish date appended)
$ws.Range("D9").Value = "1400-10-30 (2)"
$ws.Range("E9").Value = "1401-03-04 (8)"
$ws.Range("F9").Value = "1401-04-29 (2)"
$ws.Range("G9").Value = "1401-08-30 (4)"
$ws.Range("H9").Value = "1401-10-28 (2)"
$ws.Range("I9").Value = "1402-02-30 (8)"
$ws.Range("J9").Value = "1401-04-29"
$ws.Range("K9").Value = "1401-08-30 (2)"
$ws.Range("L9").Value = "1401-10-28"
$ws.Range("M9").Value = "1402-02-30 (2)"

# Row 11
$ws.Range("D11").Value = 8040580
$ws.Range("E11").Value = 10395584
$ws.Range("F11").Value = 3194446
$ws.Range("G11").Value = 11649272
$ws.Range("H11").Value = 17821794
$ws.Range("I11").Value = 21171738
$ws.Range("J11").Value = 5629545
$ws.Range("K11").Value = 11853149
$ws.Range("L11").Value = 18880390
$ws.Range("M11").Value = 25320705

# Row 12
$ws.Range("D12").Value = -3956162
$ws.Range("E12").Value = -5016199
$ws.Range("F12").Value = -1619907
$ws.Range("G12").Value = -5096967
$ws.Range("H12").Value = -8026220
$ws.Range("I12").Value = -9937618
$ws.Range("J12").Value = -3075672
$ws.Range("K12").Value = -6267993
$ws.Range("L12").Value = -10224659
$ws.Range("M12").Value = -13409448

# Row 13
$ws.Range("D13").Value = 4084418
$ws.Range("E13").Value = 5379385
$ws.Range("F13").Value = 1574539
$ws.Range("G13").Value = 6552305
$ws.Range("H13").Value = 9795574
$ws.Range("I13").Value = 11234120
$ws.Range("J13").Value = 2553873
$ws.Range("K13").Value = 5585156
$ws.Range("L13").Value = 8655731
$ws.Range("M13").Value = 11911257

# Row 14
$ws.Range("D14").Value = -309979
$ws.Range("E14").Value = -403929
$ws.Range("F14").Value = -175101
$ws.Range("G14").Value = -392745
$ws.Range("H14").Value = -626888
$ws.Range("I14").Value = -839072
$ws.Range("J14").Value = -230152
$ws.Range("K14").Value = -526967
$ws.Range("L14").Value = -779359
$ws.Range("M14").Value = -971858

# Row 15
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0

# Row 16
$ws.Range("D16").Value = -17919
$ws.Range("E16").Value = -17261
$ws.Range("F16").Value = -22480
$ws.Range("G16").Value = -24056
$ws.Range("H16").Value = -5424
$ws.Range("I16").Value = -52392
$ws.Range("J16").Value = 10856
$ws.Range("K16").Value = 2232
$ws.Range("L16").Value = -6132
$ws.Range("M16").Value = -9710

# Row 17
$ws.Range("D17").Value = 3756520
$ws.Range("E17").Value = 4958195
$ws.Range("F17").Value = 1376958
$ws.Range("G17").Value = 6135504
$ws.Range("H17").Value = 9163262
$ws.Range("I17").Value = 10342656
$ws.Range("J17").Value = 2334577
$ws.Range("K17").Value = 5060421
$ws.Range("L17").Value = 7870240
$ws.Range("M17").Value = 10929689

# Row 18
$ws.Range("D18").Value = -405105
$ws.Range("E18").Value = -531568
$ws.Range("F18").Value = -170508
$ws.Range("G18").Value = -370051
$ws.Range("H18").Value = -677718
$ws.Range("I18").Value = -959126
$ws.Range("J18").Value = -346200
$ws.Range("K18").Value = -752970
$ws.Range("L18").Value = -1166257
$ws.Range("M18").Value = -1783257

# Row 19
$ws.Range("D19").Value = 171872
$ws.Range("E19").Value = 191129
$ws.Range("F19").Value = 105196
$ws.Range("G19").Value = 31626
$ws.Range("H19").Value = 50286
$ws.Range("I19").Value = 62162
$ws.Range("J19").Value = 5528
$ws.Range("K19").Value = 196860
$ws.Range("L19").Value = 188160
$ws.Range("M19").Value = 265292

# Row 20
$ws.Range("D20").Value = 3523287
$ws.Range("E20").Value = 4617756
$ws.Range("F20").Value = 1311646
$ws.Range("G20").Value = 5797079
$ws.Range("H20").Value = 8535830
$ws.Range("I20").Value = 9445692
$ws.Range("J20").Value = 1993905
$ws.Range("K20").Value = 4504311
$ws.Range("L20").Value = 6892143
$ws.Range("M20").Value = 9411724

# Row 21
$ws.Range("D21").Value = -738127
$ws.Range("E21").Value = -642189
$ws.Range("F21").Value = -258431
$ws.Range("G21").Value = -1243936
$ws.Range("H21").Value = -1835747
$ws.Range("I21").Value = -1360326
$ws.Range("J21").Value = -356605
$ws.Range("K21").Value = -967037
$ws.Range("L21").Value = -1503365
$ws.Range("M21").Value = -1498403

# Row 22
$ws.Range("D22").Value = 2785160
$ws.Range("E22").Value = 3975567
$ws.Range("F22").Value = 1053215
$ws.Range("G22").Value = 4553143
$ws.Range("H22").Value = 6700083
$ws.Range("I22").Value = 8085366
$ws.Range("J22").Value = 1637300
$ws.Range("K22").Value = 3537274
$ws.Range("L22").Value = 5388778
$ws.Range("M22").Value = 7913321

# Row 23
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0

# Row 24
$ws.Range("D24").Value = 2785160
$ws.Range("E24").Value = 3975567
$ws.Range("F24").Value = 1053215
$ws.Range("G24").Value = 4553143
$ws.Range("H24").Value = 6700083
$ws.Range("I24").Value = 8085366
$ws.Range("J24").Value = 1637300
$ws.Range("K24").Value = 3537274
$ws.Range("L24").Value = 5388778
$ws.Range("M24").Value = 7913321

# Row 25
$ws.Range("D25").Value = 2476
$ws.Range("E25").Value = 3534
$ws.Range("F25").Value = 936
$ws.Range("G25").Value = 4047
$ws.Range("H25").Value = 5956
$ws.Range("I25").Value = 7187
$ws.Range("J25").Value = 1455
$ws.Range("K25").Value = 1572
$ws.Range("L25").Value = 2395
$ws.Range("M25").Value = 3517

# Row 26
$ws.Range("D26").Value = 1125000
$ws.Range("E26").Value = 1125000
$ws.Range("F26").Value = 1125000
$ws.Range("G26").Value = 1125000
$ws.Range("H26").Value = 1125000
$ws.Range("I26").Value = 1125000
$ws.Range("J26").Value = 1125000
$ws.Range("K26").Value = 2250000
$ws.Range("L26").Value = 2250000
$ws.Range("M26").Value = 2250000

# Row 27
$ws.Range("D27").Value = 1238
$ws.Range("E27").Value = 1767
$ws.Range("F27").Value = 468
$ws.Range("G27").Value = 2024
$ws.Range("H27").Value = 2978
$ws.Range("I27").Value = 3594
$ws.Range("J27").Value = 728
$ws.Range("K27").Value = 1572
$ws.Range("L27").Value = 2395
$ws.Range("M27").Value = 3517

